$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.420.96'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.64%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.342.82'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.17%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.12%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '189.66'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +4.62%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '563.68'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.83%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.591'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.92%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.333.61'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.17%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.186'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +1.08%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.593'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.62%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '48.04'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.93%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000274'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +3.85%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '8.73'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +1.96%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.873.98'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.15%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '604.79'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.59%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '18.18'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.59%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '66.495.72'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.64%  '

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +1.52%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.352.90'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.11%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.21'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -1.49%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.919'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +1.75%  '

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +11.97%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.18'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +3.43%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '101.06'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +1.37%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '4.03'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.31%  '

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +4.88%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.98'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.28%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.78'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +5.91%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.72'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.68%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '30.64'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.22%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.88'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +10.52%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.96'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +6.36%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '577.80'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +8.30%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '11.19'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +1.78%  '

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +2.12%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '57.34'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.64%  '

$ws.Range("B38").Value = 'Dai'
$ws.Range("C38").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.999'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.23%  '

$ws.Range("B39").Value = 'Maker'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.708.84'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -1.74%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0₃0736'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +4.19%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '34.32'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +7.94%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.131'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +5.56%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.32'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -2.50%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.72'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +2.79%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.42'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +9.06%  '

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +2.35%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0428'

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.33'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +4.24%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.130'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +1.46%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.61'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.60%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.00'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.02%  '
